$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value2 = 5.1
$ws.Range("H2").Value2 = 1.54
$ws.Range("I2").Value2 = 1.67
$ws.Range("J2").Value2 = 3.8
$ws.Range("K2").Value2 = 5
$ws.Range("M2").Value2 = 1.05
$ws.Range("N2").Value2 = 4.1
$ws.Range("O2").Value2 = 1.23
$ws.Range("P2").Value2 = 2.2
$ws.Range("Q2").Value2 = 1.67
$ws.Range("R2").Value2 = 1.48
$ws.Range("S2").Value2 = 2.68
$ws.Range("T2").Value2 = 1.77
$ws.Range("U2").Value2 = 2.04
$ws.Range("V2").Value2 = 2.5
$ws.Range("X2").Value2 = 23
$ws.Range("Y2").Value2 = 10.5
$ws.Range("Z2").Value2 = 11
$ws.Range("AB2").Value2 = 26
$ws.Range("AC2").Value2 = 11
$ws.Range("AE2").Value2 = 17
$ws.Range("AH2").Value2 = 22
$ws.Range("AI2").Value2 = 34
$ws.Range("AO2").Value2 = 8.800000000000001

# Row 4
$ws.Range("Y4").Value2 = 46
$ws.Range("AH4").Value2 = 26

# Row 5
$ws.Range("G5").Value2 = 5.4
$ws.Range("H5").Value2 = 1.74
$ws.Range("I5").Value2 = 1.76
$ws.Range("J5").Value2 = 4.1
$ws.Range("L5").Value2 = 1.31
$ws.Range("V5").Value2 = 2.32
$ws.Range("AJ5").Value2 = 120

# Row 8
$ws.Range("P8").Value2 = 2.06
$ws.Range("Q8").Value2 = 1.91
$ws.Range("R8").Value2 = 1.41
$ws.Range("AN8").Value2 = 17.5

# Row 9
$ws.Range("H9").Value2 = 9.800000000000001
$ws.Range("I9").Value2 = 10
$ws.Range("P9").Value2 = 3.65
$ws.Range("V9").Value2 = 1.11
$ws.Range("AA9").Value2 = 310
$ws.Range("AO9").Value2 = 70

# Row 10
$ws.Range("G10").Value2 = 1.4
$ws.Range("H10").Value2 = 9.800000000000001
$ws.Range("L10").Value2 = 1.26
$ws.Range("O10").Value2 = 1.17
$ws.Range("P10").Value2 = 2.72
$ws.Range("S10").Value2 = 2.36
$ws.Range("V10").Value2 = 1.1
$ws.Range("W10").Value2 = 3.5
$ws.Range("X10").Value2 = 25

# Row 11
$ws.Range("F11").Value2 = 1.68
$ws.Range("G11").Value2 = 1.69
$ws.Range("H11").Value2 = 5.4
$ws.Range("I11").Value2 = 5.6
$ws.Range("N11").Value2 = 5.8
$ws.Range("P11").Value2 = 2.62
$ws.Range("S11").Value2 = 2.46
$ws.Range("T11").Value2 = 1.65
$ws.Range("U11").Value2 = 2.46
$ws.Range("V11").Value2 = 1.21
$ws.Range("W11").Value2 = 2.44
$ws.Range("Y11").Value2 = 27
$ws.Range("AA11").Value2 = 130
$ws.Range("AE11").Value2 = 60
$ws.Range("AJ11").Value2 = 17.5
$ws.Range("AN11").Value2 = 7
$ws.Range("AO11").Value2 = 46

# Row 12
$ws.Range("F12").Value2 = 2.28
$ws.Range("G12").Value2 = 2.3
$ws.Range("J12").Value2 = 3.7
$ws.Range("K12").Value2 = 3.75
$ws.Range("L12").Value2 = 1.36
$ws.Range("N12").Value2 = 4.6
$ws.Range("P12").Value2 = 2.26
$ws.Range("Q12").Value2 = 1.77
$ws.Range("R12").Value2 = 1.48
$ws.Range("S12").Value2 = 2.96
$ws.Range("U12").Value2 = 2.42
$ws.Range("W12").Value2 = 1.76
$ws.Range("X12").Value2 = 17.5
$ws.Range("AB12").Value2 = 12.5
$ws.Range("AJ12").Value2 = 30
$ws.Range("AO12").Value2 = 27

# Row 13
$ws.Range("F13").Value2 = 7.8
$ws.Range("G13").Value2 = 8
$ws.Range("H13").Value2 = 1.46
$ws.Range("I13").Value2 = 1.47
$ws.Range("K13").Value2 = 5.4
$ws.Range("Q13").Value2 = 1.54
$ws.Range("T13").Value2 = 1.77
$ws.Range("U13").Value2 = 2.22
$ws.Range("AA13").Value2 = 13.5
$ws.Range("AF13").Value2 = 70
$ws.Range("AG13").Value2 = 27
$ws.Range("AK13").Value2 = 90
$ws.Range("AM13").Value2 = 90

# Row 14
$ws.Range("G14").Value2 = 2.78
$ws.Range("Q14").Value2 = 1.69
$ws.Range("R14").Value2 = 1.58
$ws.Range("S14").Value2 = 2.7
$ws.Range("W14").Value2 = 1.56

# Row 16
$ws.Range("F16").Value2 = 1.89
$ws.Range("G16").Value2 = 2.02
$ws.Range("Q16").Value2 = 1.79
$ws.Range("T16").Value2 = 1.74
$ws.Range("U16").Value2 = 2.2
$ws.Range("W16").Value2 = 1.98

# Row 18
$ws.Range("F18").Value2 = 1.8

# Row 19
$ws.Range("F19").Value2 = 2.1
$ws.Range("G19").Value2 = 2.26
$ws.Range("H19").Value2 = 4.6
$ws.Range("I19").Value2 = 5.2
$ws.Range("J19").Value2 = 2.9
$ws.Range("K19").Value2 = 3.15
$ws.Range("M19").Value2 = 1.15
$ws.Range("N19").Value2 = 2.36
$ws.Range("P19").Value2 = 1.45
$ws.Range("V19").Value2 = 1.23
$ws.Range("W19").Value2 = 1.79
$ws.Range("X19").Value2 = 7.6
$ws.Range("AB19").Value2 = 6.4
$ws.Range("AC19").Value2 = 7.6
$ws.Range("AF19").Value2 = 980
$ws.Range("AN19").Value2 = 42

# Row 20
$ws.Range("L20").Value2 = 1.4

# Row 21
$ws.Range("F21").Value2 = 1.62
$ws.Range("G21").Value2 = 1.74
$ws.Range("H21").Value2 = 5.3
$ws.Range("I21").Value2 = 6.6
$ws.Range("K21").Value2 = 5.2
$ws.Range("L21").Value2 = 1.34
$ws.Range("O21").Value2 = 1.26
$ws.Range("P21").Value2 = 2.1
$ws.Range("Q21").Value2 = 1.65
$ws.Range("R21").Value2 = 1.43
$ws.Range("S21").Value2 = 2.78
$ws.Range("T21").Value2 = 1.78
$ws.Range("U21").Value2 = 2.04
$ws.Range("V21").Value2 = 1.19
$ws.Range("W21").Value2 = 2.34
$ws.Range("AB21").Value2 = 1000
$ws.Range("AF21").Value2 = 1000
$ws.Range("AG21").Value2 = 1000
$ws.Range("AM21").Value2 = 120
